$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range first so numeric-looking
# strings (prices like "242.00", hour codes like "20") are stored
# exactly as typed instead of being normalised into numbers -- this
# matches the original workbook where every data cell is inline text.
$ws.Range("B2:G51").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = 'BNB'
$ws.Range("C2").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D2").Value = '242.00'
$ws.Range("E2").Value = '1BNBBNB'
$ws.Range("F2").Value = '25-12-2022'
$ws.Range("G2").Value = '20'

# Row 3
$ws.Range("B3").Value = 'OKB'
$ws.Range("C3").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D3").Value = '22.86'
$ws.Range("E3").Value = '2OKBOKB'
$ws.Range("F3").Value = '25-12-2022'
$ws.Range("G3").Value = '20'

# Row 4
$ws.Range("B4").Value = 'LEO'
$ws.Range("C4").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D4").Value = '3.611'
$ws.Range("E4").Value = '3LEOLEO'
$ws.Range("F4").Value = '25-12-2022'
$ws.Range("G4").Value = '20'

# Row 5
$ws.Range("B5").Value = 'HuobiToken'
$ws.Range("C5").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D5").Value = '5.379'
$ws.Range("E5").Value = '4HuobiTokenHT'
$ws.Range("F5").Value = '25-12-2022'
$ws.Range("G5").Value = '20'

# Row 6
$ws.Range("B6").Value = 'Cronos'
$ws.Range("C6").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D6").Value = '0.05938'
$ws.Range("E6").Value = '5CronosCRO'
$ws.Range("F6").Value = '25-12-2022'
$ws.Range("G6").Value = '20'

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '3.402'
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("F7").Value = '25-12-2022'
$ws.Range("G7").Value = '20'

# Row 8
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '6.444'
$ws.Range("E8").Value = '7KuCoinTokenKCS'
$ws.Range("F8").Value = '25-12-2022'
$ws.Range("G8").Value = '20'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.8047'
$ws.Range("E9").Value = '8MXTokenMX'
$ws.Range("F9").Value = '25-12-2022'
$ws.Range("G9").Value = '20'

# Row 10
$ws.Range("B10").Value = 'FTXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D10").Value = '0.9144'
$ws.Range("E10").Value = '9FTXTokenFTT'
$ws.Range("F10").Value = '25-12-2022'
$ws.Range("G10").Value = '20'

# Row 11
$ws.Range("B11").Value = 'One'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D11").Value = '0.01111'
$ws.Range("E11").Value = '10OneONE'
$ws.Range("F11").Value = '25-12-2022'
$ws.Range("G11").Value = '20'

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = '0.1410'
$ws.Range("E12").Value = '11WazirXWRX'
$ws.Range("F12").Value = '25-12-2022'
$ws.Range("G12").Value = '20'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '0.07422'
$ws.Range("E13").Value = '12MandalaExchangeTokenMDX'
$ws.Range("F13").Value = '25-12-2022'
$ws.Range("G13").Value = '20'

# Row 14
$ws.Range("B14").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C14").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D14").Value = '0.03265'
$ws.Range("E14").Value = '13LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("F14").Value = '25-12-2022'
$ws.Range("G14").Value = '20'

# Row 15
$ws.Range("B15").Value = 'BitrueCoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D15").Value = '0.03034'
$ws.Range("E15").Value = '14BitrueCoinBTR'
$ws.Range("F15").Value = '25-12-2022'
$ws.Range("G15").Value = '20'

# Row 16
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = '0.09311'
$ws.Range("E16").Value = '15BitMartTokenBMX'
$ws.Range("F16").Value = '25-12-2022'
$ws.Range("G16").Value = '20'

# Row 17
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D17").Value = '3.871'
$ws.Range("E17").Value = '16MCDexMCB'
$ws.Range("F17").Value = '25-12-2022'
$ws.Range("G17").Value = '20'

# Row 18
$ws.Range("B18").Value = 'BitForexToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D18").Value = '0.001579'
$ws.Range("E18").Value = '17BitForexTokenBF'
$ws.Range("F18").Value = '25-12-2022'
$ws.Range("G18").Value = '20'

# Row 19
$ws.Range("B19").Value = 'CoinExToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D19").Value = '0.04481'
$ws.Range("E19").Value = '18CoinExTokenCET'
$ws.Range("F19").Value = '25-12-2022'
$ws.Range("G19").Value = '20'

# Row 20
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '0.006100'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("F20").Value = '25-12-2022'
$ws.Range("G20").Value = '20'

# Row 21
$ws.Range("B21").Value = 'UpBots'
$ws.Range("C21").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D21").Value = '0.007493'
$ws.Range("E21").Value = '20UpBotsUBXTBestin24h'
$ws.Range("F21").Value = '25-12-2022'
$ws.Range("G21").Value = '20'

# Row 22
$ws.Range("B22").Value = 'HotbitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D22").Value = '0.004411'
$ws.Range("E22").Value = '21HotbitTokenHTB'
$ws.Range("F22").Value = '25-12-2022'
$ws.Range("G22").Value = '20'

# Row 23
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").Value = '0.0009808'
$ws.Range("E23").Value = '22BitKanKAN'
$ws.Range("F23").Value = '25-12-2022'
$ws.Range("G23").Value = '20'

# Row 24
$ws.Range("B24").Value = 'NitroEx'
$ws.Range("C24").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D24").Value = '0.00007801'
$ws.Range("E24").Value = '23NitroExNTX'
$ws.Range("F24").Value = '25-12-2022'
$ws.Range("G24").Value = '20'

# Row 25
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.137'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("F25").Value = '25-12-2022'
$ws.Range("G25").Value = '20'

# Row 26
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '0.3248'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("F26").Value = '25-12-2022'
$ws.Range("G26").Value = '20'

# Row 27
$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").Value = '0.1297'
$ws.Range("E27").Value = '26ProBitTokenPROB'
$ws.Range("F27").Value = '25-12-2022'
$ws.Range("G27").Value = '20'

# Row 28
$ws.Range("B28").Value = 'Spectre.aiUtilityToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("D28").Value = '--'
$ws.Range("E28").Value = '27Spectre.aiUtilityTokenSXUT'
$ws.Range("F28").Value = '25-12-2022'
$ws.Range("G28").Value = '20'

# Row 29
$ws.Range("B29").Value = 'LegolasExchange'
$ws.Range("C29").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("D29").Value = '--'
$ws.Range("E29").Value = '28LegolasExchangeLGO'
$ws.Range("F29").Value = '25-12-2022'
$ws.Range("G29").Value = '20'

# Row 30
$ws.Range("B30").Value = 'BitZToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("D30").Value = '--'
$ws.Range("E30").Value = '29BitZTokenBZ'
$ws.Range("F30").Value = '25-12-2022'
$ws.Range("G30").Value = '20'

# Row 31
$ws.Range("B31").Value = 'Birake'
$ws.Range("C31").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("D31").Value = '--'
$ws.Range("E31").Value = '30BirakeBIR'
$ws.Range("F31").Value = '25-12-2022'
$ws.Range("G31").Value = '20'

# Row 32
$ws.Range("B32").Value = 'ZBToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D32").Value = '--'
$ws.Range("E32").Value = '31ZBTokenZB'
$ws.Range("F32").Value = '25-12-2022'
$ws.Range("G32").Value = '20'

# Row 33
$ws.Range("B33").Value = 'NashExchange'
$ws.Range("C33").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range("D33").Value = '--'
$ws.Range("E33").Value = '32NashExchangeNEX'
$ws.Range("F33").Value = '25-12-2022'
$ws.Range("G33").Value = '20'

# Row 34
$ws.Range("B34").Value = 'AAXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("D34").Value = '--'
$ws.Range("E34").Value = '33AAXTokenAAB'
$ws.Range("F34").Value = '25-12-2022'
$ws.Range("G34").Value = '20'

# Row 35
$ws.Range("B35").Value = 'CenX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range("D35").Value = '--'
$ws.Range("E35").Value = '34CenXCENX'
$ws.Range("F35").Value = '25-12-2022'
$ws.Range("G35").Value = '20'

# Row 36
$ws.Range("B36").Value = 'BNIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range("D36").Value = '--'
$ws.Range("E36").Value = '35BNIXTokenBNIX'
$ws.Range("F36").Value = '25-12-2022'
$ws.Range("G36").Value = '20'

# Row 37
$ws.Range("B37").Value = 'Polkally'
$ws.Range("C37").Value = 'https://coinranking.com/coin/NkDWUL8F-+polkally-kally'
$ws.Range("D37").Value = '--'
$ws.Range("E37").Value = '36PolkallyKALLY'
$ws.Range("F37").Value = '25-12-2022'
$ws.Range("G37").Value = '20'

# Row 38
$ws.Range("B38").Value = 'Charli3'
$ws.Range("C38").Value = 'https://coinranking.com/coin/8SgjMSqUk+charli3-c3'
$ws.Range("D38").Value = '--'
$ws.Range("E38").Value = '37Charli3C3'
$ws.Range("F38").Value = '25-12-2022'
$ws.Range("G38").Value = '20'

# Row 39
$ws.Range("B39").Value = 'BlubitexToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Y9oImHIW5+blubitextoken-bbe'
$ws.Range("D39").Value = '--'
$ws.Range("E39").Value = '38BlubitexTokenBBE'
$ws.Range("F39").Value = '25-12-2022'
$ws.Range("G39").Value = '20'

# Row 40
$ws.Range("B40").Value = 'IDEX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ZiRElvGxqQaf+idex-idex'
$ws.Range("D40").Value = '0.03859'
$ws.Range("E40").Value = '39IDEXIDEX'
$ws.Range("F40").Value = '25-12-2022'
$ws.Range("G40").Value = '20'

# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.006127'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("F41").Value = '25-12-2022'
$ws.Range("G41").Value = '20'

# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1064'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("F42").Value = '25-12-2022'
$ws.Range("G42").Value = '20'

# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '0.002801'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("F43").Value = '25-12-2022'
$ws.Range("G43").Value = '20'

# Row 44
$ws.Range("B44").Value = 'LocalTraders'
$ws.Range("C44").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D44").Value = '0.007226'
$ws.Range("E44").Value = '43LocalTradersLCT'
$ws.Range("F44").Value = '25-12-2022'
$ws.Range("G44").Value = '20'

# Row 45
$ws.Range("B45").Value = 'CoinLion'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion'
$ws.Range("D45").Value = '0.00005194'
$ws.Range("E45").Value = '44CoinLionLION'
$ws.Range("F45").Value = '25-12-2022'
$ws.Range("G45").Value = '20'

# Row 46
$ws.Range("B46").Value = 'Kangarootoken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar'
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").Value = '45KangarootokenGAR'
$ws.Range("F46").Value = '25-12-2022'
$ws.Range("G46").Value = '20'

# Row 47
$ws.Range("B47").Value = 'ACDXExchange'
$ws.Range("C47").Value = 'https://coinranking.com/coin/-y35lbZ7U+acdxexchange-acxt'
$ws.Range("D47").Value = '0.0005802'
$ws.Range("E47").Value = '46ACDXExchangeACXT'
$ws.Range("F47").Value = '25-12-2022'
$ws.Range("G47").Value = '20'

# Row 48
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '0.9582'
$ws.Range("E48").Value = '47CoinbaseStockTokenCOIN'
$ws.Range("F48").Value = '25-12-2022'
$ws.Range("G48").Value = '20'

# Row 49
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '0.002264'
$ws.Range("E49").Value = '48BOLOBOLO'
$ws.Range("F49").Value = '25-12-2022'
$ws.Range("G49").Value = '20'

# Row 50
$ws.Range("B50").Value = 'CryptobidCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '49CryptobidCoinCBC'
$ws.Range("F50").Value = '25-12-2022'
$ws.Range("G50").Value = '20'

# Row 51
$ws.Range("B51").Value = 'SpecialPowerGold'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg'
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '50SpecialPowerGoldSPGWorstin24h'
$ws.Range("F51").Value = '25-12-2022'
$ws.Range("G51").Value = '20'

